$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Package info": update version/build metadata, add a new "Name" row
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")

$wsInfo.Cells.Item(4, 2).Value = "1.2.0"
$wsInfo.Cells.Item(5, 2).Value = "2.34.4"

$wsInfo.Cells.Item(6, 1).Value = "DHIS2 build"
$wsInfo.Cells.Item(6, 2).Value = "aff07fb"

$wsInfo.Cells.Item(7, 1).Value = "Last updated"
$wsInfo.Cells.Item(7, 2).Value = "20210520T090044"

$wsInfo.Cells.Item(8, 1).Value = "Name"
$wsInfo.Cells.Item(8, 2).Value = "MAL-HIST_CUSTOM_V1.2.0_2.34.4-en"

# column B is narrower now
$wsInfo.Columns.Item(2).ColumnWidth = 34.7109375

# ---------------------------------------------------------------------------
# Sheet "dataElements": rows reshuffled/renamed/recoded (the "Last updated"
# column F and "Categorycombo" column E are identical on every row, both
# before and after, so they are left untouched)
# ---------------------------------------------------------------------------
$wsDE = $wb.Worksheets.Item("dataElements")

$deRows = @(
    @("MAL - Plasmodium falciparum (Mic+RDT)", "P. falciparum (micr+RDT)", "MAL_PF_MICR_RDT", "Cases confirmed as P.falciparum positive with microscopy and/or RDT", "IIU1O0Z4l49"),
    @("MAL - Mixed/Other malaria species (Mic+RDT)", "Other species (micr+RDT)", "MAL_MIX_OTHER_SPECIES_MICR_RDT", "Cases confirmed as P.malariae or P. ovale or P. knowlesi with microscopy and/or RDT", "JkOyLRb3dpX"),
    @("MAL - Plasmodium vivax (Mic+RDT)", "P. vivax (micr+RDT)", "MAL_PV_MICR_RDT", "Cases confirmed as P.vivax positive with microscopy", "pUC8tgzn0lV"),
    @("MAL - Mixed malaria species  (Mic+RDT)", "Mixed (micr+RDT)", "MAL_MIX_SPECIES_MICR_RDT", "Cases confirmed as mixed infection with P.falciparum and P.vivax with microscopy and/or RDT", "TNTW2ruEVEu"),
    @("MAL - Malaria tested cases (Mic+RDT)", "Tested (micr+RDT)", "MAL_TEST_CASES_MICR_RDT", "Suspected cases tested with both microscopy and/or RDT", "tuOTgWfDO6m"),
    @("MAL - Malaria confirmed cases (Mic+RDT)", "Positive (micr+RDT)", "MAL_CONFI_CASES_MICR_RDT", "Cases confirmed as positive with microscopy and/or RDT", "X0luAFiy268")
)

for ($i = 0; $i -lt $deRows.Count; $i++) {
    $r = $i + 2
    $row = $deRows[$i]
    $wsDE.Cells.Item($r, 1).Value = $row[0]
    $wsDE.Cells.Item($r, 2).Value = $row[1]
    $wsDE.Cells.Item($r, 3).Value = $row[2]
    $wsDE.Cells.Item($r, 4).Value = $row[3]
    $wsDE.Cells.Item($r, 7).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet "dataElementGroups": reorder data element names to match the new order
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")

$wsDEG.Cells.Item(2, 2).Value = "MAL - Plasmodium falciparum (Mic+RDT)"
$wsDEG.Cells.Item(3, 2).Value = "MAL - Mixed/Other malaria species (Mic+RDT)"
$wsDEG.Cells.Item(4, 2).Value = "MAL - Plasmodium vivax (Mic+RDT)"
$wsDEG.Cells.Item(5, 2).Value = "MAL - Mixed malaria species  (Mic+RDT)"
$wsDEG.Cells.Item(6, 2).Value = "MAL - Malaria tested cases (Mic+RDT)"
$wsDEG.Cells.Item(7, 2).Value = "MAL - Malaria confirmed cases (Mic+RDT)"

# ---------------------------------------------------------------------------
# Sheet "userGroups": reorder rows 2/4 and refresh "Last updated" dates.
# A leading single-quote forces the date-looking string to stay text
# (matches the original t="str" cells) instead of being auto-converted
# into an Excel date serial number.
# ---------------------------------------------------------------------------
$wsUG = $wb.Worksheets.Item("userGroups")

$ugRows = @(
    @("Malaria data capture", "'2021-05-20", "fRSrUJ6SMGH"),
    @("Malaria admin", "'2021-05-20", "suMb19wGXPR"),
    @("Malaria access", "'2021-05-20", "ZXEVDM9XRea")
)

for ($i = 0; $i -lt $ugRows.Count; $i++) {
    $r = $i + 2
    $row = $ugRows[$i]
    $wsUG.Cells.Item($r, 1).Value = $row[0]
    $wsUG.Cells.Item($r, 2).Value = $row[1]
    $wsUG.Cells.Item($r, 3).Value = $row[2]
}
